$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Consolidated "Absent" column (H) values formed from the attendance data.
$ws.Range("H3").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H10").Value = 1
$ws.Range("H14").Value = 0
